$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix data sheets back to root s of 63 GeV: column F ("rs") rows 2-18 were 125, now 63
$ws.Range("F2:F18").Value = 63

# Touch row 25 (formatting only, no visible value) so the sheet's used range
# extends down to row 25, matching the new dimension A1:J25
$ws.Cells.Item(25, 1).NumberFormat = "General"

# Move the active cell / selection to G23
$ws.Range("G23").Select()
